# Insert a new data row at row 237 (pushes existing rows 237:249 down to 238:250)
# and populate it with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(237).Insert()

$ws.Range("A237").Value = 9
$ws.Range("B237").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C237").Value = "Metropolitana"
$ws.Range("D237").Value = 44753
$ws.Range("E237").Value = 13
$ws.Range("F237").Value = 100112003
$ws.Range("G237").Value = "Ajo"
$ws.Range("H237").Value = "Chino"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 430
$ws.Range("K237").Value = 22000
$ws.Range("L237").Value = 24000
$ws.Range("M237").Value = 23000
$ws.Range("N237").Value = '$/caja 10 kilos'
$ws.Range("O237").Value = "China"
$ws.Range("P237").Value = 2300
$ws.Range("Q237").Value = 10
$ws.Range("R237").Value = "Hortaliza"
